# Re-pull data / push all data / mean calculation
# Update the dSF (column F) values for the specified rows to reflect
# the re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 0
    5  = -3
    8  = -2
    11 = -2
    20 = -2
    23 = 0
    27 = -5
    30 = 4
    31 = -1
    38 = 0
    41 = 2
    46 = -2
    47 = 1
    48 = 3
    49 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
